# Moved functionality to strategy
# Add two new backlog rows (describing strategy-related work items) and
# refresh the existing "pending" status AutoFilter so it matches the
# current result set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11: "Deleting an event has the appropriate behavior if calendar is on day view"
# Write the Name (column B) first so its shared string precedes the "?" priority
# marker in the shared-strings table, matching the authoring order.
$ws.Range("B11").Value = "Deleting an event has the appropriate behavior if calendar is on day view"
$ws.Range("A11").Value = "?"
$ws.Range("C11").Value = "pending"

# New row 12: "Either the work day is 8 hours, or the 9 = 100% (fix the status bar)"
$ws.Range("B12").Value = "Either the work day is 8 hours, or the 9 = 100% (fix the status bar)"
$ws.Range("A12").Value = "?"
$ws.Range("C12").Value = "pending"

# Reapply the existing Status AutoFilter (column C, 3rd column of A1:C9) so the
# "pending" criteria is recorded as a standard values-filter instead of the
# stale blank-inclusive filter list.
$ws.Range("A1:C9").AutoFilter(3, @("pending"), 7)

# Re-applying the filter recalculates row visibility for the whole filtered
# range; row 9 ("current") was left visible in the original workbook despite
# not matching "pending", so restore that pre-existing (intentional) state.
$ws.Rows.Item(9).Hidden = $false

# Leave the selection where the author left off after entering the new rows.
$ws.Range("B14").Select()
